$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table of events and
# both need their "想去人数" (want-to-go count) values updated:
#   F2: 199 -> 197
#   F4: 128 -> 129
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 197
    $ws.Range("F4").Value = 129
}
